$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slight precision drift on the existing row 9 timestamp
$ws.Range("A9").Value = 45862.75021494213

# Append the new row 10 recorded by the scheduled task
$ws.Range("A10").Value = 45862.79190448848
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B10").Value = 2025
$ws.Range("C10").Value = 30
$ws.Range("D10").Value = 19.33
$ws.Range("E10").Value = 74.44
$ws.Range("F10").Value = 87.89
$ws.Range("G10").Value = 13.54
$ws.Range("H10").Value = "ESE"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "19:00:20"
